$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to text format so numeric-looking strings
# such as "1.004" are stored as text, matching the source data which
# uses inline strings rather than numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "23.578.37"
$ws.Range("E2").Value = "  +1.44%  "
$ws.Range("D3").Value = "1.646.36"
$ws.Range("E3").Value = "  +2.55%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.38%  "
$ws.Range("D5").Value = "1.003"
$ws.Range("E5").Value = "  +0.19%  "
$ws.Range("D6").Value = "305.42"
$ws.Range("E6").Value = "  +0.30%  "
$ws.Range("D7").Value = "0.3776"
$ws.Range("E7").Value = "  +0.37%  "
$ws.Range("D8").Value = "52.94"
$ws.Range("E8").Value = "  +0.92%  "
$ws.Range("D9").Value = "0.3676"
$ws.Range("E9").Value = "  +1.49%  "
$ws.Range("D10").Value = "1.267"
$ws.Range("E10").Value = "  -0.57%  "
$ws.Range("D11").Value = "0.08157"
$ws.Range("E11").Value = "  +0.05%  "
$ws.Range("D13").Value = "23.13"
$ws.Range("E13").Value = "  +0.76%  "
$ws.Range("D14").Value = "6.713"
$ws.Range("E14").Value = "  +1.61%  "
$ws.Range("D15").Value = "0.00001271"
$ws.Range("E15").Value = "  +1.80%  "
$ws.Range("D16").Value = "7.391"
$ws.Range("E16").Value = "  +0.17%  "
$ws.Range("D17").Value = "1.650.01"
$ws.Range("E17").Value = "  +2.94%  "
$ws.Range("D18").Value = "95.10"
$ws.Range("E18").Value = "  +1.10%  "
$ws.Range("D19").Value = "0.06920"
$ws.Range("E19").Value = "  -0.07%  "
$ws.Range("D20").Value = "18.36"
$ws.Range("E20").Value = "  +1.01%  "
$ws.Range("D21").Value = "6.602"
$ws.Range("E21").Value = "  +0.86%  "
$ws.Range("D22").Value = "1.003"
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D23").Value = "23.603.01"
$ws.Range("E23").Value = "  +1.59%  "
$ws.Range("D24").Value = "12.97"
$ws.Range("E24").Value = "  +0.29%  "
$ws.Range("D25").Value = "3.255"
$ws.Range("E25").Value = "  +6.16%  "
$ws.Range("D26").Value = "2.431"
$ws.Range("E26").Value = "  -0.74%  "
$ws.Range("D27").Value = "21.51"
$ws.Range("E27").Value = "  +1.53%  "
$ws.Range("D28").Value = "152.14"
$ws.Range("E28").Value = "  +1.11%  "
$ws.Range("D29").Value = "5.326"
$ws.Range("E29").Value = "  +0.96%  "
$ws.Range("D30").Value = "137.54"
$ws.Range("E30").Value = "  +1.69%  "
$ws.Range("D31").Value = "2.311"
$ws.Range("E31").Value = "  -3.37%  "
$ws.Range("D32").Value = "7.035"
$ws.Range("E32").Value = "  +4.46%  "
$ws.Range("D33").Value = "1.833.23"
$ws.Range("E33").Value = "  +3.03%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "0.9783"
$ws.Range("E34").Value = "  +1.66%  "
$ws.Range("B35").Value = "FraxShare"
$ws.Range("C35").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D35").Value = "11.01"
$ws.Range("E35").Value = "  +5.99%  "
$ws.Range("D36").Value = "0.02891"
$ws.Range("E36").Value = "  +4.20%  "
$ws.Range("D37").Value = "6.371"
$ws.Range("E37").Value = "  +4.10%  "
$ws.Range("D38").Value = "0.2585"
$ws.Range("E38").Value = "  +2.48%  "
$ws.Range("D39").Value = "0.07332"
$ws.Range("E39").Value = "  -2.09%  "
$ws.Range("D40").Value = "0.08880"
$ws.Range("E40").Value = "  +0.85%  "
$ws.Range("D41").Value = "1.386"
$ws.Range("E41").Value = "  -1.72%  "
$ws.Range("D42").Value = "0.7217"
$ws.Range("E42").Value = "  +1.67%  "
$ws.Range("D43").Value = "12.76"
$ws.Range("E43").Value = "  +2.42%  "
$ws.Range("D44").Value = "16.58"
$ws.Range("E44").Value = "  +3.92%  "
$ws.Range("D45").Value = "0.6652"
$ws.Range("E45").Value = "  +1.67%  "
$ws.Range("D46").Value = "2.392"
$ws.Range("E46").Value = "  +2.37%  "
$ws.Range("D47").Value = "1.002"
$ws.Range("E47").Value = "  +0.22%  "
$ws.Range("D48").Value = "4.022"
$ws.Range("E48").Value = "  +0.34%  "
$ws.Range("D49").Value = "0.08053"
$ws.Range("E49").Value = "  +1.28%  "
$ws.Range("D50").Value = "1.230"
$ws.Range("E50").Value = "  +2.07%  "
$ws.Range("D51").Value = "128.98"
$ws.Range("E51").Value = "  -3.77%  "
